# Right-align every cell paragraph in the first table ("Evolution de
# l'encours des ressources ...") that is currently left-aligned.
#
# Note: $d.Tables.Item(1).Range.Paragraphs does not reliably scope to the
# table's own paragraphs in this host, so we walk the table via
# Rows/Columns + Cell(r,c).Range.Paragraphs instead, which is scoped
# correctly.

$d = $word.ActiveDocument

$table = $d.Tables.Item(1)

$wdAlignParagraphLeft = 0
$wdAlignParagraphRight = 2

$changed = 0

for ($r = 1; $r -le $table.Rows.Count; $r++) {
    for ($c = 1; $c -le $table.Columns.Count; $c++) {
        $cell = $table.Cell($r, $c)
        $para = $cell.Range.Paragraphs.Item(1)
        if ($para.Alignment -eq $wdAlignParagraphLeft) {
            $para.Alignment = $wdAlignParagraphRight
            $changed = $changed + 1
        }
    }
}

Write-Output "Right-aligned $changed paragraph(s) in table 1."
